$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header row (e.g. A1) for the new headers
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous
$headerRange.Borders.Weight = 2           # xlThin

# Fill season record values for each data row (rows 2-52)
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 57   # AD
    $ws.Cells.Item($r, 31).Value = 105  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
